$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Slide 16 - "Interpretive Compilers"
#   "Example - Oracle/Sun Java Development Kit" -> "Example - Java Development Kit"
# ----------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(4)
$tr16 = $sh16.TextFrame.TextRange
$s16para2 = $tr16.Paragraphs(2, 1)
$s16para2.Text = "Example – Java Development Kit"

# ----------------------------------------------------------------------
# Slide 17 - "Just-In-Time Compiler"
# ----------------------------------------------------------------------
$s17 = $p.Slides.Item(17)

# Body placeholder ("Rectangle 3")
$sh17body = $s17.Shapes.Item(4)

# Nudge the placeholder position/size (small manual adjustment in the source deck)
$sh17body.Left   = 36.12504007007874
$sh17body.Top    = 107.37504197007874
$sh17body.Width  = 648.0
$sh17body.Height = 388.62496062992125

$tr17body = $sh17body.TextFrame.TextRange

# Paragraph 1, 3rd run: "... is a compiler that converts ... just before the program is run."
#                     -> "... is a compiler that converts ... as the program is running."
$run3 = $tr17body.Characters(30, 105)
$run3.Text = " is a compiler that converts program source code into native machine code as the program is running."

# Paragraph 3 (lvl 1): "Use of the JIT compiler is optional."
#                    -> "The JVM interpreter starts executing initially with no delay."
$para3 = $tr17body.Paragraphs(3, 1)
$para3.Text = "The JVM interpreter starts executing initially with no delay."

# Paragraph 4 (lvl 1): "Translation for a method is performed when the method is first called."
#                    -> "Methods that are executed frequently (hot) are JIT compiled."
$tr17body_b = $sh17body.TextFrame.TextRange
$para4 = $tr17body_b.Paragraphs(4, 1)
$para4.Text = "Methods that are executed frequently (hot) are JIT compiled."

# New paragraph (lvl 1) inserted right after paragraph 4:
#   "Execution switches to the compiled version once it becomes available,"
$tr17body_c = $sh17body.TextFrame.TextRange
$para4b = $tr17body_c.Paragraphs(4, 1)
$para4b.InsertAfter("`rExecution switches to the compiled version once it becomes available,")

# TextBox 1 (the boxed note below)
$sh17note = $s17.Shapes.Item(5)

$sh17note.Left   = 58.64055118110236
$sh17note.Top    = 437.83748031496066
$sh17note.Width  = 602.7189763779528
$sh17note.Height = 58.16251968503937

$tr17note = $sh17note.TextFrame.TextRange
$note1 = $tr17note.Paragraphs(1, 1)
$note1.Font.Size = 21
$tr17note_b = $sh17note.TextFrame.TextRange
$note2 = $tr17note_b.Paragraphs(2, 1)
$note2.Font.Size = 21

# ----------------------------------------------------------------------
# Slide 9 - "Integrated Development Environment (IDE)"
#   Merge "Apache " + "Netbeans" + ", and Microsoft Visual Studio."
#   into a single run "Apache NetBeans, and Microsoft Visual Studio."
# ----------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(5)
$tr9 = $sh9.TextFrame.TextRange
$s9para1 = $tr9.Paragraphs(1, 1)

# Remove paragraph 2 (and the paragraph mark that separates it from paragraph 1)
$s9toDelete = $tr9.Characters($s9para1.Start + $s9para1.Length, 46)
$s9toDelete.Delete()

# Re-insert a clean single-run paragraph with the corrected text
$tr9_b = $sh9.TextFrame.TextRange
$s9para1b = $tr9_b.Paragraphs(1, 1)
$s9para1b.InsertAfter("`rApache NetBeans, and Microsoft Visual Studio.")
